# Daily attendance processing - 2025-10-25 13:45:11
# Normalize the "Recorded By" (column G) cell values: for a small, fixed set of
# "recorder list" strings, the tokens have been reordered (e.g. "System" moved
# ahead of specific user emails). This applies that exact, deterministic
# re-mapping to every matching cell in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact old-value -> new-value replacements observed for the "Recorded By" column.
$map = @{
    "system, backup@backdoor.com, System" = "system, System, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
